$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values with new latency data
$ws.Range("A2").Value = 0.5496434448710916
$ws.Range("B2").Value = 1689839664.873628
$ws.Range("C2").Value = 0.9767541235291549
$ws.Range("D2").Value = 1689839664.992314
$ws.Range("E2").Value = 0.1186857223510742

# Delete rows 3 through 13 (entire rows)
$ws.Range("A3:E13").EntireRow.Delete()
